# Refactor synthetic array: insert a "statut_name" column right after
# "statut_label" (i.e. before the old "NCTId" column), shifting the
# NCTId..intervention_type columns one place to the right, and fill the
# new column with a human readable status label derived from the
# results / results_3y flags of each clinical trial row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row so the script isn't hard-coded to 19 rows.
$lastRow = $ws.Cells(1, 1).SpecialCells(11).Row

# Insert a new column before column C (NCTId), pushing NCTId..intervention_type
# from C..L to D..M.
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "statut_name"

# After the insert, "results_3y" now lives in column K and "results" in
# column L (they used to be J and K respectively before the insert).
for ($row = 2; $row -le $lastRow; $row++) {
    $results3y = $ws.Cells.Item($row, 11).Value2
    $results = $ws.Cells.Item($row, 12).Value2

    if ($results) {
        if ($results3y) {
            $statusText = "résultat et / ou publication posté dans les 36 mois"
        } else {
            $statusText = "résultat et / ou publication posté"
        }
    } else {
        $statusText = "pas de résultat ni de publication"
    }

    $ws.Cells.Item($row, 3).Value = $statusText
}
